# ICTU-Kwaliteitsaanpak.pptx edit script
#
# 1. Bump the version string on the title slide (slide 1).
# 2. Reword the M02 slide (slide 8): title becomes "bewaakt continu dat" and
#    the body text is tightened up accordingly.
# 3. Remove the M06 slide (slide 20) entirely — its content now overlaps too
#    much with the reworded M02, per the commit message.

$p = $ppt.ActivePresentation

# --- 1. Title slide version bump -----------------------------------------
$slide1 = $p.Slides.Item(1)
$slide1.Shapes.Item(2).TextFrame.TextRange.Text = "Versie 4.0.0-dev, 06-06-2023"

# --- 2. M02 slide: update title + body ------------------------------------
$slide8 = $p.Slides.Item(8)
$slide8.Shapes.Item(1).TextFrame.TextRange.Text = "M02: Het project bewaakt continu dat het product aan de kwaliteitsnormen voldoet"
$slide8.Shapes.Item(2).TextFrame.TextRange.Text = "Projecten bewaken zo snel mogelijk vanaf de start de door het project en ICTU vastgestelde kwaliteitsnormen en voldoen daar zo snel en goed mogelijk aan. De kwaliteit van producten, die nog niet zijn afgerond of nog niet aan de normen voldoen, wordt door het project bewaakt. Het voldoen aan de kwaliteitsnormen is onderdeel van de Definition of Done en herstel van de kwaliteit wordt planmatig opgepakt."

# --- 3. Remove the M06 slide -----------------------------------------------
$slide20 = $p.Slides.Item(20)
$slide20.Delete()
